$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.802.60"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.909.75"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'312.76"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.5172"
$ws.Range("E7").Value = "  +5.51%  "
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").Value = "'0.07255"
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").Value = "'21.28"
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("D11").Value = "'0.9055"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07658"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.934.45"
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("D14").Value = "'5.457"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "'92.10"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D17").Value = "'0.000008714"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "'0.9996"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "27.837.15"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").Value = "'5.156"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").Value = "2.185.30"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").Value = "'6.646"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "'154.04"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("D27").Value = "'2.174"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("D28").Value = "'18.37"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "'115.02"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'4.863"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "'0.09080"
$ws.Range("E31").Value = "  +1.76%  "
$ws.Range("D32").Value = "'3.188"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "'4.853"
$ws.Range("E33").Value = "  +5.04%  "
$ws.Range("D34").Value = "'1.234"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "'0.7811"
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("D36").Value = "'0.02097"
$ws.Range("E36").Value = "  +2.84%  "
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("D39").Value = "'0.5587"
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("D40").Value = "'1.095"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").Value = "'0.05294"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").Value = "'6.732"
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("D43").Value = "'115.38"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").Value = "'8.568"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "'0.1519"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "'0.4828"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").Value = "'10.50"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").Value = "'0.9997"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("D50").Value = "'67.02"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "'0.05994"
$ws.Range("E51").Value = "  -1.02%  "
